$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.123.91"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.607.17"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").Value = "'212.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Value = "'0.483"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").Value = "'0.0620"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "'18.43"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.19%  "
$ws.Range("D11").Value = "'0.0792"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "1.830.18"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "1.592.99"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "'0.512"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "26.105.91"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'60.87"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").Value = "'1.00"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").Value = "'197.89"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.68%  "
$ws.Range("D21").Value = "'4.26"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("E24").Value = "  +2.54%  "
$ws.Range("D25").Value = "'142.65"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "'15.24"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'1.17"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("D31").Value = "'0.0476"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("D33").Value = "'3.04"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("E34").Value = "  +4.31%  "
$ws.Range("D35").Value = "'2.34"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").Value = "1.105.36"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "'2.37"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("D39").Value = "'0.508"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").Value = "'0.789"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("E42").Value = "  +8.52%  "
$ws.Range("D43").Value = "'5.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("D44").Value = "1.743.17"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "'93.16"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("E46").Value = "  -5.27%  "
$ws.Range("D47").Value = "'1.55"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.53%  "
$ws.Range("D48").Value = "'53.86"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D49").Value = "'0.0509"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "'0.407"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -0.76%  "
